$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E28) is reversed in order (2001..2101 -> 2101..2001)
$ws.Range("E16").Value = "2101"
$ws.Range("E17").Value = "2012"
$ws.Range("E18").Value = "2011"
$ws.Range("E19").Value = "2010"
$ws.Range("E20").Value = "2009"
$ws.Range("E21").Value = "2008"
$ws.Range("E22").Value = "2007"
$ws.Range("E23").Value = "2006"
$ws.Range("E24").Value = "2005"
$ws.Range("E25").Value = "2004"
$ws.Range("E26").Value = "2003"
$ws.Range("E27").Value = "2002"
$ws.Range("E28").Value = "2001"

# "Valor Mora" column (F16:F28) keeps the same values but the one different
# amount now sits on the first row instead of the last
$ws.Range("F16").Value = 25013
$ws.Range("F28").Value = 31266
